$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Total"
$ws.Range("C29").Formula = "=SUM(C2:C26)"

$ws.Range("A31").Value = "NYPD proportion"
$ws.Range("C31").Formula = "=C3/C29"
$ws.Range("C31").NumberFormat = "0.0000000000"

$ws.Columns.Item(2).ColumnWidth = 9.875
$ws.Columns.Item(3).ColumnWidth = 17.25

$ws.Range("C32").Select()
